$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 28 (swap)
$ws.Range("B28").Value = "Cosmos"
$ws.Range("C28").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").Value = "'6.760"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.05%  "

# Row 29 (swap)
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").Value = "'15.56"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.13%  "

# Row 2
$ws.Range("D2").Value = "'25.902.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.44%  "

# Row 3
$ws.Range("D3").Value = "'1.632.07"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.50%  "

# Row 4
$ws.Range("E4").Value = "  -0.13%  "

# Row 5
$ws.Range("D5").Value = "'215.95"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.61%  "

# Row 6
$ws.Range("D6").Value = "'0.5116"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.52%  "

# Row 7
$ws.Range("E7").Value = "  -0.02%  "

# Row 8
$ws.Range("D8").Value = "'0.2566"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.08%  "

# Row 9
$ws.Range("D9").Value = "'0.06342"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.30%  "

# Row 10
$ws.Range("D10").Value = "'19.46"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -0.19%  "

# Row 11
$ws.Range("D11").Value = "'0.07774"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.02%  "

# Row 12
$ws.Range("D12").Value = "'4.256"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.70%  "

# Row 13
$ws.Range("D13").Value = "'1.632.84"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -0.67%  "

# Row 14
$ws.Range("D14").Value = "'1.857.09"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.49%  "

# Row 15
$ws.Range("D15").Value = "'0.5510"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.22%  "

# Row 16
$ws.Range("D16").Value = "'63.74"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.70%  "

# Row 17
$ws.Range("D17").Value = "'0.0₅7631"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -1.24%  "

# Row 18
$ws.Range("D18").Value = "'25.918.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.50%  "

# Row 20
$ws.Range("D20").Value = "'194.82"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.10%  "

# Row 21
$ws.Range("D21").Value = "'4.418"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -0.18%  "

# Row 22
$ws.Range("D22").Value = "'9.858"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.71%  "

# Row 23
$ws.Range("D23").Value = "'6.022"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "

# Row 24
$ws.Range("E24").Value = "  -0.20%  "

# Row 25
$ws.Range("E25").Value = "  +0.80%  "

# Row 26
$ws.Range("D26").Value = "'142.16"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.70%  "

# Row 27
$ws.Range("D27").Value = "'0.1255"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +5.05%  "

# Row 30
$ws.Range("E30").Value = "  +0.59%  "

# Row 31
$ws.Range("D31").Value = "'0.04914"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +0.87%  "

# Row 32
$ws.Range("D32").Value = "'3.229"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.71%  "

# Row 33
$ws.Range("D33").Value = "'3.184"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +0.48%  "

# Row 34
$ws.Range("D34").Value = "'1.544"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.14%  "

# Row 35
$ws.Range("D35").Value = "'2.371"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.48%  "

# Row 36
$ws.Range("D36").Value = "'0.8978"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.19%  "

# Row 37
$ws.Range("D37").Value = "'0.5509"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.97%  "

# Row 38
$ws.Range("D38").Value = "'2.535"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -1.50%  "

# Row 39
$ws.Range("D39").Value = "'1.113.40"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.41%  "

# Row 40
$ws.Range("D40").Value = "'0.01557"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -0.19%  "

# Row 41
$ws.Range("D41").Value = "'1.000"
$ws.Range("D41").Style = "Normal"

# Row 42
$ws.Range("D42").Value = "'5.575"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +3.01%  "

# Row 43
$ws.Range("D43").Value = "'0.7974"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.49%  "

# Row 44
$ws.Range("D44").Value = "'97.64"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.48%  "

# Row 45
$ws.Range("D45").Value = "'1.766.22"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.62%  "

# Row 46
$ws.Range("E46").Value = "  -8.38%  "

# Row 47
$ws.Range("E47").Value = "  -2.03%  "

# Row 48
$ws.Range("E48").Value = "  +0.20%  "

# Row 49
$ws.Range("D49").Value = "'54.74"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.30%  "

# Row 50
$ws.Range("E50").Value = "  +1.45%  "

# Row 51
$ws.Range("D51").Value = "'7.578"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +3.46%  "

Write-Host "Done applying changes"